$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'thermal men leggings'
$ws.Cells.Item(2, 1).Value = 'compression pants with knee pads'
$ws.Cells.Item(3, 1).Value = 'under armor compression pants men'
$ws.Cells.Item(4, 1).Value = 'compression pants with pads'
$ws.Cells.Item(5, 1).Value = 'basketball warm up pants men'
$ws.Cells.Item(6, 1).Value = 'basketball compression leggings'
$ws.Cells.Item(7, 1).Value = 'bjj compression pants'
$ws.Cells.Item(8, 1).Value = 'padded leggings basketball'
$ws.Cells.Item(9, 1).Value = 'football tights'
$ws.Cells.Item(10, 1).Value = 'knee pads for volleyball youth'
$ws.Cells.Item(11, 1).Value = 'recovery compression pants men'
$ws.Cells.Item(12, 1).Value = 'hockey compression leggings'
$ws.Cells.Item(13, 1).Value = 'knee compressions'
$ws.Cells.Item(14, 1).Value = 'men sports tights'
$ws.Cells.Item(15, 1).Value = 'wrestling knee pads men'
$ws.Cells.Item(16, 1).Value = 'mens running leggings'
$ws.Cells.Item(17, 1).Value = 'thick knee pad'
$ws.Cells.Item(18, 1).Value = 'football pads pants'
$ws.Cells.Item(19, 1).Value = 'boys athletic pants'
$ws.Cells.Item(20, 1).Value = 'running compression tights'
$ws.Cells.Item(21, 1).Value = 'volleyball gear men'
$ws.Cells.Item(22, 1).Value = 'girls capri leggings'
$ws.Cells.Item(23, 1).Value = 'black capri leggings'
$ws.Cells.Item(24, 1).Value = 'coolomg basketball knee pads'
$ws.Cells.Item(25, 1).Value = 'tight with knee pads'
$ws.Cells.Item(26, 1).Value = 'tights mens nike'
$ws.Cells.Item(27, 1).Value = 'muscle compression pants for men'
$ws.Cells.Item(28, 1).Value = 'hayabusa compression pants men'
$ws.Cells.Item(29, 1).Value = 'girls athletic leggings'
$ws.Cells.Item(30, 1).Value = 'black capri leggings for women'
$ws.Cells.Item(31, 1).Value = 'asics youth knee pads'
$ws.Cells.Item(32, 1).Value = 'men''s tights sports'
$ws.Cells.Item(33, 1).Value = 'tights with knee pads'
$ws.Cells.Item(34, 1).Value = 'lavento compression pants'
$ws.Cells.Item(35, 1).Value = 'tsla compression pants men'
$ws.Cells.Item(36, 1).Value = 'mens black baseball pants'
$ws.Cells.Item(37, 1).Value = 'knee pad protector'
$ws.Cells.Item(38, 1).Value = 'compression pants boys basketball'
$ws.Cells.Item(39, 1).Value = 'black leggings men'
$ws.Cells.Item(40, 1).Value = 'athletic leggings men'
$ws.Cells.Item(41, 1).Value = 'youth hex knee pads'
$ws.Cells.Item(42, 1).Value = 'wrestling youth knee pads'
$ws.Cells.Item(43, 1).Value = 'adult volleyball knee pads'
$ws.Cells.Item(44, 1).Value = 'basketball shorts with pads'
$ws.Cells.Item(45, 1).Value = 'mens wrestling tights'
$ws.Cells.Item(46, 1).Value = 'softball sliding pants'
$ws.Cells.Item(47, 1).Value = 'black baseball pants'
$ws.Cells.Item(48, 1).Value = 'wrestling tights for men'
$ws.Cells.Item(49, 1).Value = 'mens running capris'
$ws.Cells.Item(50, 1).Value = 'mens football pants with pads'
$ws.Cells.Item(51, 1).Value = 'soccer pants'
$ws.Cells.Item(52, 1).Value = 'baseball sliding shorts boys'
$ws.Cells.Item(53, 1).Value = 'compression knee guards'
$ws.Cells.Item(54, 1).Value = 'knees protector'
$ws.Cells.Item(55, 1).Value = 'baseball pants youth'
$ws.Cells.Item(56, 1).Value = 'mens basketball shorts black'
$ws.Cells.Item(57, 1).Value = 'mens compression'
$ws.Cells.Item(58, 1).Value = 'wrestling clothes for men'
$ws.Cells.Item(59, 1).Value = 'mens sports pants'
$ws.Cells.Item(60, 1).Value = 'boys gym pants'
$ws.Cells.Item(61, 1).Value = 'knee pads football adult'
$ws.Cells.Item(62, 1).Value = 'leggings men short'
$ws.Cells.Item(63, 1).Value = 'running pants men tall'
$ws.Cells.Item(64, 1).Value = 'boys running tights youth'
$ws.Cells.Item(65, 1).Value = 'softball pants for girls youth'
$ws.Cells.Item(66, 1).Value = 'little boys athletic pants'
$ws.Cells.Item(67, 1).Value = 'boys running pants size'
$ws.Cells.Item(68, 1).Value = 'football youth pants'
$ws.Cells.Item(69, 1).Value = 'comfortable knee pads'
$ws.Cells.Item(70, 1).Value = 'knee sleeve with pad'
$ws.Cells.Item(71, 1).Value = 'women compression tights'
$ws.Cells.Item(72, 1).Value = 'knee pads toddler'
$ws.Cells.Item(73, 1).Value = 'knee pads basketball kids'
$ws.Cells.Item(74, 1).Value = 'knee protector for kids'
$ws.Cells.Item(75, 1).Value = 'dark green knee pads for basketball'
$ws.Cells.Item(76, 1).Value = 'skateboarding knee pads youth'
$ws.Cells.Item(77, 1).Value = 'compression pants with padding basketball'
$ws.Cells.Item(78, 1).Value = 'firefighter compression pants'
$ws.Cells.Item(79, 1).Value = 'skateboard knee and elbow pads youth'
$ws.Cells.Item(80, 1).Value = 'skateboard youth knee pads'
$ws.Cells.Item(81, 1).Value = 'warm up pants men'
$ws.Cells.Item(82, 1).Value = 'mc david knee pad'
$ws.Cells.Item(83, 1).Value = 'compression knee pads men'
$ws.Cells.Item(84, 1).Value = 'men basketball pants'
$ws.Cells.Item(85, 1).Value = 'premium knee pad'
$ws.Cells.Item(86, 1).Value = 'kids compression pants for basketball'
$ws.Cells.Item(87, 1).Value = 'toddler knee pad'
$ws.Cells.Item(88, 1).Value = 'padded knee sleeves men'
$ws.Cells.Item(89, 1).Value = 'mtb knee pads for men'
$ws.Cells.Item(90, 1).Value = 'compression knee sleeves pads'
$ws.Cells.Item(91, 1).Value = 'basketball knee pads leggings'
$ws.Cells.Item(92, 1).Value = 'soccer compression pants men'
$ws.Cells.Item(93, 1).Value = 'mens small leggings'
$ws.Cells.Item(94, 1).Value = 'compressions pants mens'
$ws.Cells.Item(95, 1).Value = 'capri basketball leggings for boys'
$ws.Cells.Item(96, 1).Value = 'youth small black baseball pants'
$ws.Cells.Item(97, 1).Value = 'mens black leggings'
$ws.Cells.Item(98, 1).Value = 'basket ball knee pads youth'
$ws.Cells.Item(99, 1).Value = 'bjj leggings men'
$ws.Cells.Item(100, 1).Value = 'mens gym tights'
